# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el texto de conversión del día con las nuevas tasas ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.47 = 39345.64 pesos`n✅ 39345.64 pesos = 9.43 = 970.68 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Hoja "tasas": actualizar tasas N10/O10 y N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 105.6
$ws2.Range("O10").Value = 4154.9

$ws2.Range("N12").Value = 4170.95
$ws2.Range("O12").Value = 102.9
